$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.319.27'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.691.20'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '217.83'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5400'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").Value = '0.06448'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = '21.68'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = '0.07668'
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = '1.690.87'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '4.539'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '0.5795'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008386'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").Value = '67.09'
$ws.Range("E16").Value = '  +3.52%  '
$ws.Range("D17").Value = '26.346.27'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.910'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '10.87'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '190.35'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '6.261'
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '148.95'
$ws.Range("E24").Value = '  +2.28%  '
$ws.Range("D25").Value = '0.1289'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").Value = '7.868'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '15.89'
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("D28").Value = '0.06287'
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("D29").Value = '1.371'
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").Value = '1.327'
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").Value = '3.602'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '3.585'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = '1.677'
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").Value = '1.033'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("D36").Value = '2.417'
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").Value = '2.765'
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").Value = '0.01655'
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").Value = '1.111.08'
$ws.Range("D40").Value = '6.122'
$ws.Range("E40").Value = '  -5.19%  '
$ws.Range("D41").Value = '0.8822'
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("E45").Value = '  -3.14%  '
$ws.Range("D46").Value = '57.66'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.142'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").Value = '6.052'
$ws.Range("E51").Value = '  -0.47%  '
